$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Shooting und Bewegung Verknüpfen"
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3

$ws.Range("B18").Select()
